# ---------------------------------------------------------------------------
# Proposed Automation Testing Approach.pptx - apply commit changes
#   * refresh the "last saved" date field (2/8/2019 -> 6/9/2019) on the
#     slide master and every slide layout
#   * remove the stray empty "Text Placeholder" shapes (idx 11/13/14, and
#     the empty duplicate Title) that were left behind on slides 3-9
#   * reposition the two screenshot pictures that moved on slides 6 and 8
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Date placeholder text: 2/8/2019 -> 6/9/2019 -----------------------
$newDate = "6/9/2019"

$design = $p.Designs.Item(1)
$master = $design.SlideMaster

for ($j = 1; $j -le $master.Shapes.Count; $j++) {
    $sh = $master.Shapes.Item($j)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $sh = $layout.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- 2) Remove stray empty placeholder shapes ------------------------------
function Remove-ShapeByIdName {
    param($slide, [int]$shapeId, [string]$shapeName)
    for ($k = $slide.Shapes.Count; $k -ge 1; $k--) {
        $cand = $slide.Shapes.Item($k)
        if ($cand.Id -eq $shapeId -and $cand.Name -eq $shapeName) {
            $cand.Delete()
            break
        }
    }
}

Remove-ShapeByIdName -slide $p.Slides.Item(3) -shapeId 7 -shapeName "Text Placeholder 6"
Remove-ShapeByIdName -slide $p.Slides.Item(4) -shapeId 6 -shapeName "Text Placeholder 5"
Remove-ShapeByIdName -slide $p.Slides.Item(5) -shapeId 7 -shapeName "Text Placeholder 6"
Remove-ShapeByIdName -slide $p.Slides.Item(6) -shapeId 7 -shapeName "Text Placeholder 6"
Remove-ShapeByIdName -slide $p.Slides.Item(7) -shapeId 7 -shapeName "Text Placeholder 6"
Remove-ShapeByIdName -slide $p.Slides.Item(8) -shapeId 3 -shapeName "Text Placeholder 2"
Remove-ShapeByIdName -slide $p.Slides.Item(8) -shapeId 5 -shapeName "Text Placeholder 4"
Remove-ShapeByIdName -slide $p.Slides.Item(9) -shapeId 3 -shapeName "Text Placeholder 2"
Remove-ShapeByIdName -slide $p.Slides.Item(9) -shapeId 4 -shapeName "Title 3"
Remove-ShapeByIdName -slide $p.Slides.Item(9) -shapeId 5 -shapeName "Text Placeholder 4"

# --- 3) Reposition moved pictures ------------------------------------------
# Values below are expressed in points (EMU / 12700) but nudged by a few
# millionths so that the COM layer's single-precision round trip lands on
# the exact target EMU instead of one EMU short.

function Set-ShapePositionEmu {
    param($slide, [int]$shapeId, [string]$shapeName, [double]$leftPt, [double]$topPt)
    for ($k = 1; $k -le $slide.Shapes.Count; $k++) {
        $cand = $slide.Shapes.Item($k)
        if ($cand.Id -eq $shapeId -and $cand.Name -eq $shapeName) {
            $cand.Left = $leftPt
            $cand.Top = $topPt
            break
        }
    }
}

# Slide 6 - Picture 7: off x=5934008,y=1645918 -> x=6656005,y=1724776
Set-ShapePositionEmu -slide $p.Slides.Item(6) -shapeId 8 -shapeName "Picture 7" -leftPt 524.0948818897638 -topPt 135.80913585826772

# Slide 8 - Picture 6: off x=7577703,y=2449880 -> x=7577703,y=1616161
Set-ShapePositionEmu -slide $p.Slides.Item(8) -shapeId 7 -shapeName "Picture 6" -leftPt 596.6695275590552 -topPt 127.25677165354331
